$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, pushing existing rows 94-128 down to 95-129
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new price-report record
$ws.Cells.Item(94, 1).Value = 4
$ws.Cells.Item(94, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value = "Los Lagos"
$ws.Cells.Item(94, 4).Value = 44988
$ws.Cells.Item(94, 5).Value = 10
$ws.Cells.Item(94, 6).Value = 100112031
$ws.Cells.Item(94, 7).Value = "Poroto verde"
$ws.Cells.Item(94, 8).Value = "Magnum"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 40
$ws.Cells.Item(94, 11).Value = 30000
$ws.Cells.Item(94, 12).Value = 30000
$ws.Cells.Item(94, 13).Value = 30000
$ws.Cells.Item(94, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(94, 15).Value = "Región Metropolitana"
$ws.Cells.Item(94, 16).Value = 1200
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"
